$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TPM-derived metrics for rows 2-21 (Mfng-Notch2 LR pairs)
# Values correspond to the new TPM-based recalculation.

# Row 2
$ws.Range("G2").Value = 12.227772
$ws.Range("H2").Value = 36.683316
$ws.Range("I2").Value = 0.6223179025846677
$ws.Range("J2").Value = 0.6245980896688198
$ws.Range("M2").Value = 1.400501333333333
$ws.Range("N2").Value = 4.201504
$ws.Range("O2").Value = 0.00926314904242919
$ws.Range("P2").Value = 0.009687730200823723
$ws.Range("Q2").Value = 17.125010989696
$ws.Range("R2").Value = 154.125098907264
$ws.Range("S2").Value = 0.005764623483413707
$ws.Range("T2").Value = 0.006050937776661429

# Row 3
$ws.Range("G3").Value = 12.227772
$ws.Range("H3").Value = 36.683316
$ws.Range("I3").Value = 0.6223179025846677
$ws.Range("J3").Value = 0.6245980896688198
$ws.Range("O3").Value = 0.1405812059498714
$ws.Range("P3").Value = 0.1470248171880475
$ws.Range("Q3").Value = 259.895925868092
$ws.Range("R3").Value = 2339.063332812828
$ws.Range("S3").Value = 0.08748620122954717
$ws.Range("T3").Value = 0.09183141994956191

# Row 4
$ws.Range("G4").Value = 12.227772
$ws.Range("H4").Value = 36.683316
$ws.Range("I4").Value = 0.6223179025846677
$ws.Range("J4").Value = 0.6245980896688198
$ws.Range("M4").Value = 63.87756733333333
$ws.Range("N4").Value = 191.632702
$ws.Range("O4").Value = 0.4224968677952986
$ws.Range("P4").Value = 0.4418622271050682
$ws.Range("Q4").Value = 781.0803292666481
$ws.Range("R4").Value = 7029.722963399833
$ws.Range("S4").Value = 0.2629273646149619
$ws.Range("T4").Value = 0.2759863029466358

# Row 5
$ws.Range("G5").Value = 12.227772
$ws.Range("H5").Value = 36.683316
$ws.Range("I5").Value = 0.6223179025846677
$ws.Range("J5").Value = 0.6245980896688198
$ws.Range("M5").Value = 19.878555
$ws.Range("N5").Value = 39.75711
$ws.Range("O5").Value = 0.1314800731212866
$ws.Range("P5").Value = 0.0916710195312133
$ws.Range("Q5").Value = 243.07043822946
$ws.Range("R5").Value = 1458.42262937676
$ws.Range("S5").Value = 0.0818224033365178
$ws.Range("T5").Value = 0.0572575436771889

# Row 6
$ws.Range("G6").Value = 12.227772
$ws.Range("H6").Value = 36.683316
$ws.Range("I6").Value = 0.6223179025846677
$ws.Range("J6").Value = 0.6245980896688198
$ws.Range("M6").Value = 44.77944466666667
$ws.Range("N6").Value = 134.338334
$ws.Range("O6").Value = 0.2961787040911142
$ws.Range("P6").Value = 0.3097542059748472
$ws.Range("Q6").Value = 547.5528396706161
$ws.Range("R6").Value = 4927.975557035545
$ws.Range("S6").Value = 0.1843173099202272
$ws.Range("T6").Value = 0.1934718853187717

# Row 7
$ws.Range("G7").Value = 3.888411
$ws.Range("H7").Value = 11.665233
$ws.Range("I7").Value = 0.1978960499023984
$ws.Range("J7").Value = 0.1986211455731449
$ws.Range("M7").Value = 1.400501333333333
$ws.Range("N7").Value = 4.201504
$ws.Range("O7").Value = 0.00926314904242919
$ws.Range("P7").Value = 0.009687730200823723
$ws.Range("Q7").Value = 5.445724790048
$ws.Range("R7").Value = 49.01152311043199
$ws.Range("S7").Value = 0.001833140605153921
$ws.Range("T7").Value = 0.001924188070491161

# Row 8
$ws.Range("G8").Value = 3.888411
$ws.Range("H8").Value = 11.665233
$ws.Range("I8").Value = 0.1978960499023984
$ws.Range("J8").Value = 0.1986211455731449
$ws.Range("O8").Value = 0.1405812059498714
$ws.Range("P8").Value = 0.1470248171880475
$ws.Range("Q8").Value = 82.64646879257099
$ws.Range("R8").Value = 743.8182191331389
$ws.Range("S8").Value = 0.02782046534799509
$ws.Range("T8").Value = 0.02920223761757219

# Row 9
$ws.Range("G9").Value = 3.888411
$ws.Range("H9").Value = 11.665233
$ws.Range("I9").Value = 0.1978960499023984
$ws.Range("J9").Value = 0.1986211455731449
$ws.Range("M9").Value = 63.87756733333333
$ws.Range("N9").Value = 191.632702
$ws.Range("O9").Value = 0.4224968677952986
$ws.Range("P9").Value = 0.4418622271050682
$ws.Range("Q9").Value = 248.382235472174
$ws.Range("R9").Value = 2235.440119249566
$ws.Range("S9").Value = 0.08361046123282544
$ws.Range("T9").Value = 0.08776318173310975

# Row 10
$ws.Range("G10").Value = 3.888411
$ws.Range("H10").Value = 11.665233
$ws.Range("I10").Value = 0.1978960499023984
$ws.Range("J10").Value = 0.1986211455731449
$ws.Range("M10").Value = 19.878555
$ws.Range("N10").Value = 39.75711
$ws.Range("O10").Value = 0.1314800731212866
$ws.Range("P10").Value = 0.0916710195312133
$ws.Range("Q10").Value = 77.29599192610499
$ws.Range("R10").Value = 463.7759515566299
$ws.Range("S10").Value = 0.02601938711158112
$ws.Range("T10").Value = 0.01820780291514772

# Row 11
$ws.Range("G11").Value = 3.888411
$ws.Range("H11").Value = 11.665233
$ws.Range("I11").Value = 0.1978960499023984
$ws.Range("J11").Value = 0.1986211455731449
$ws.Range("M11").Value = 44.77944466666667
$ws.Range("N11").Value = 134.338334
$ws.Range("O11").Value = 0.2961787040911142
$ws.Range("P11").Value = 0.3097542059748472
$ws.Range("Q11").Value = 174.120885215758
$ws.Range("R11").Value = 1567.087966941822
$ws.Range("S11").Value = 0.05861259560484283
$ws.Range("T11").Value = 0.06152373523682401

# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.215192
$ws.Range("H12").Value = 0.430384
$ws.Range("I12").Value = 0.01095194072092608
$ws.Range("J12").Value = 0.007328045922130521
$ws.Range("M12").Value = 1.400501333333333
$ws.Range("N12").Value = 4.201504
$ws.Range("O12").Value = 0.00926314904242919
$ws.Range("P12").Value = 0.009687730200823723
$ws.Range("Q12").Value = 0.3013766829226667
$ws.Range("R12").Value = 1.808260097536
$ws.Range("S12").Value = 0.0001014494592017877
$ws.Range("T12").Value = 0.00007099213179284698

# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.215192
$ws.Range("H13").Value = 0.430384
$ws.Range("I13").Value = 0.01095194072092608
$ws.Range("J13").Value = 0.007328045922130521
$ws.Range("O13").Value = 0.1405812059498714
$ws.Range("P13").Value = 0.1470248171880475
$ws.Range("Q13").Value = 4.573811490712
$ws.Range("R13").Value = 27.442868944272
$ws.Range("S13").Value = 0.001539637034039293
$ws.Range("T13").Value = 0.001077404612046857

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.5
$ws.Range("G14").Value = 0.215192
$ws.Range("H14").Value = 0.430384
$ws.Range("I14").Value = 0.01095194072092608
$ws.Range("J14").Value = 0.007328045922130521
$ws.Range("M14").Value = 63.87756733333333
$ws.Range("N14").Value = 191.632702
$ws.Range("O14").Value = 0.4224968677952986
$ws.Range("P14").Value = 0.4418622271050682
$ws.Range("Q14").Value = 13.74594146959467
$ws.Range("R14").Value = 82.475648817568
$ws.Range("S14").Value = 0.004627160650871055
$ws.Range("T14").Value = 0.003237986691480805

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.5
$ws.Range("G15").Value = 0.215192
$ws.Range("H15").Value = 0.430384
$ws.Range("I15").Value = 0.01095194072092608
$ws.Range("J15").Value = 0.007328045922130521
$ws.Range("M15").Value = 19.878555
$ws.Range("N15").Value = 39.75711
$ws.Range("O15").Value = 0.1314800731212866
$ws.Range("P15").Value = 0.0916710195312133
$ws.Range("Q15").Value = 4.27770600756
$ws.Range("R15").Value = 17.11082403024
$ws.Range("S15").Value = 0.001439961966807357
$ws.Range("T15").Value = 0.000671769440853255

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.5
$ws.Range("G16").Value = 0.215192
$ws.Range("H16").Value = 0.430384
$ws.Range("I16").Value = 0.01095194072092608
$ws.Range("J16").Value = 0.007328045922130521
$ws.Range("M16").Value = 44.77944466666667
$ws.Range("N16").Value = 134.338334
$ws.Range("O16").Value = 0.2961787040911142
$ws.Range("P16").Value = 0.3097542059748472
$ws.Range("Q16").Value = 9.636178256709334
$ws.Range("R16").Value = 57.817069540256
$ws.Range("S16").Value = 0.003243731610006591
$ws.Range("T16").Value = 0.002269893045956756

# Row 17
$ws.Range("G17").Value = 3.31738
$ws.Range("H17").Value = 9.95214
$ws.Range("I17").Value = 0.168834106792008
$ws.Range("J17").Value = 0.1694527188359048
$ws.Range("M17").Value = 1.400501333333333
$ws.Range("N17").Value = 4.201504
$ws.Range("O17").Value = 0.00926314904242919
$ws.Range("P17").Value = 0.009687730200823723
$ws.Range("Q17").Value = 4.645995113173333
$ws.Range("R17").Value = 41.81395601856
$ws.Range("S17").Value = 0.001563935494659776
$ws.Range("T17").Value = 0.001641612221878286

# Row 18
$ws.Range("G18").Value = 3.31738
$ws.Range("H18").Value = 9.95214
$ws.Range("I18").Value = 0.168834106792008
$ws.Range("J18").Value = 0.1694527188359048
$ws.Range("O18").Value = 0.1405812059498714
$ws.Range("P18").Value = 0.1470248171880475
$ws.Range("Q18").Value = 70.50945557017999
$ws.Range("R18").Value = 634.58510013162
$ws.Range("S18").Value = 0.02373490233828985
$ws.Range("T18").Value = 0.02491375500886651

# Row 19
$ws.Range("G19").Value = 3.31738
$ws.Range("H19").Value = 9.95214
$ws.Range("I19").Value = 0.168834106792008
$ws.Range("J19").Value = 0.1694527188359048
$ws.Range("M19").Value = 63.87756733333333
$ws.Range("N19").Value = 191.632702
$ws.Range("O19").Value = 0.4224968677952986
$ws.Range("P19").Value = 0.4418622271050682
$ws.Range("Q19").Value = 211.9061643202533
$ws.Range("R19").Value = 1907.15547888228
$ws.Range("S19").Value = 0.07133188129664032
$ws.Range("T19").Value = 0.07487475573384184

# Row 20
$ws.Range("G20").Value = 3.31738
$ws.Range("H20").Value = 9.95214
$ws.Range("I20").Value = 0.168834106792008
$ws.Range("J20").Value = 0.1694527188359048
$ws.Range("M20").Value = 19.878555
$ws.Range("N20").Value = 39.75711
$ws.Range("O20").Value = 0.1314800731212866
$ws.Range("P20").Value = 0.0916710195312133
$ws.Range("Q20").Value = 65.94472078589999
$ws.Range("R20").Value = 395.6683247154
$ws.Range("S20").Value = 0.02219832070638031
$ws.Range("T20").Value = 0.01553390349802342

# Row 21
$ws.Range("G21").Value = 3.31738
$ws.Range("H21").Value = 9.95214
$ws.Range("I21").Value = 0.168834106792008
$ws.Range("J21").Value = 0.1694527188359048
$ws.Range("M21").Value = 44.77944466666667
$ws.Range("N21").Value = 134.338334
$ws.Range("O21").Value = 0.2961787040911142
$ws.Range("P21").Value = 0.3097542059748472
$ws.Range("Q21").Value = 148.5504341483067
$ws.Range("R21").Value = 1336.95390733476
$ws.Range("S21").Value = 0.05000506695603771
$ws.Range("T21").Value = 0.05248869237329472
